$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '63.484.26'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +1.41%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.411.89'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +2.18%  '
$ws.Range('E4').Value = '  +0.01%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '568.07'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.44%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '155.99'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +2.68%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '3.411.99'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +2.03%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.545'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('E11').Value = '  +3.93%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.431'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -0.53%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '4.000.71'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +2.18%  '
$ws.Range('E14').Value = '  -3.00%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.0000193'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +8.27%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '27.16'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +1.11%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '63.561.43'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +1.49%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '3.410.93'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +1.43%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '6.25'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -1.19%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '14.10'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +2.33%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '380.03'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.90%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '8.06'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -3.49%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +2.48%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '0.529'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('E26').Value = '  +28.25%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '9.37'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +5.14%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '0.178'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +0.22%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.18%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '6.04'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +8.58%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.37'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +5.20%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '2.00'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +1.34%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '23.25'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +1.76%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '6.37'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -2.43%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '6.79'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +1.52%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '159.76'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  -1.41%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.963.21'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +6.02%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.0764'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +3.61%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '27.17'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('E43').Value = '  +1.32%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '41.88'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +2.76%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '4.32'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +1.93%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '23.27'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +6.32%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '1.07'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +3.46%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '2.21'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +23.94%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.833'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +4.45%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '6.35'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +0.88%  '
